$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44433
$ws.Range("I2").Value = 'Segunda'
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 17000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 17500
$ws.Range("P2").Value = 972

# Row 3
$ws.Range("D3").Value = 44433
$ws.Range("I3").Value = 'Tercera'
$ws.Range("J3").Value = 120
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 14500
$ws.Range("P3").Value = 806

# Row 5
$ws.Range("D5").Value = 44377
$ws.Range("J5").Value = 100
$ws.Range("M5").Value = 17600
$ws.Range("P5").Value = 978

# Row 6
$ws.Range("D6").Value = 44533
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 6000
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = 6500
$ws.Range("P6").Value = 650

# Row 7
$ws.Range("D7").Value = 44533
$ws.Range("I7").Value = 'Segunda'
$ws.Range("J7").Value = 120
$ws.Range("K7").Value = 4000
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = 4500
$ws.Range("P7").Value = 450

# Row 8
$ws.Range("D8").Value = 44554
$ws.Range("H8").Value = 'Cultivar XV región'
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 5000
$ws.Range("L8").Value = 6000
$ws.Range("M8").Value = 5500
$ws.Range("N8").Value = '$/caja 10 kilos'
$ws.Range("O8").Value = 'Región de Arica y Parinacota'
$ws.Range("P8").Value = 550
$ws.Range("Q8").Value = 10

# Row 9
$ws.Range("D9").Value = 44391
$ws.Range("I9").Value = 'Segunda'
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 15000
$ws.Range("L9").Value = 16000
$ws.Range("M9").Value = 15500
$ws.Range("P9").Value = 861

# Row 10
$ws.Range("D10").Value = 44412
$ws.Range("H10").Value = 'Cultivar IV Región'
$ws.Range("I10").Value = 'Primera'
$ws.Range("J10").Value = 150
$ws.Range("K10").Value = 17000
$ws.Range("L10").Value = 18000
$ws.Range("M10").Value = 17500
$ws.Range("N10").Value = '$/bandeja 18 kilos'
$ws.Range("O10").Value = 'Provincia de Limarí'
$ws.Range("P10").Value = 972
$ws.Range("Q10").Value = 18

# Row 11
$ws.Range("D11").Value = 44526
$ws.Range("H11").Value = 'Cultivar XV región'
$ws.Range("K11").Value = 5000
$ws.Range("L11").Value = 5500
$ws.Range("M11").Value = 5250
$ws.Range("N11").Value = '$/caja 10 kilos'
$ws.Range("O11").Value = 'Región de Arica y Parinacota'
$ws.Range("P11").Value = 525
$ws.Range("Q11").Value = 10

# Row 12
$ws.Range("I12").Value = 'Segunda'
$ws.Range("K12").Value = 4000
$ws.Range("L12").Value = 4500
$ws.Range("M12").Value = 4250
$ws.Range("P12").Value = 425

# Row 13
$ws.Range("I13").Value = 'Tercera'
$ws.Range("J13").Value = 120
$ws.Range("K13").Value = 3000
$ws.Range("L13").Value = 3500
$ws.Range("M13").Value = 3250
$ws.Range("P13").Value = 325

# Row 14
$ws.Range("D14").Value = 44363
$ws.Range("H14").Value = 'Cultivar IV Región'
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 140
$ws.Range("K14").Value = 14000
$ws.Range("L14").Value = 15000
$ws.Range("M14").Value = 14500
$ws.Range("N14").Value = '$/bandeja 18 kilos'
$ws.Range("O14").Value = 'Provincia de Limarí'
$ws.Range("P14").Value = 806
$ws.Range("Q14").Value = 18

# Row 17
$ws.Range("D17").Value = 44211
$ws.Range("H17").Value = 'Cultivar XV región'
$ws.Range("I17").Value = 'Segunda'
$ws.Range("J17").Value = 140
$ws.Range("K17").Value = 4500
$ws.Range("L17").Value = 5000
$ws.Range("M17").Value = 4750
$ws.Range("N17").Value = '$/caja 10 kilos'
$ws.Range("O17").Value = 'Región de Arica y Parinacota'
$ws.Range("P17").Value = 475
$ws.Range("Q17").Value = 10

# Row 18
$ws.Range("D18").Value = 44405
$ws.Range("J18").Value = 140
$ws.Range("K18").Value = 17000
$ws.Range("L18").Value = 18000
$ws.Range("M18").Value = 17500
$ws.Range("P18").Value = 972

# Row 19
$ws.Range("D19").Value = 44398
$ws.Range("J19").Value = 100

# Row 20
$ws.Range("D20").Value = 44398
$ws.Range("H20").Value = 'Cultivar IV Región'
$ws.Range("I20").Value = 'Segunda'
$ws.Range("K20").Value = 15000
$ws.Range("L20").Value = 16000
$ws.Range("M20").Value = 15500
$ws.Range("N20").Value = '$/bandeja 18 kilos'
$ws.Range("O20").Value = 'Provincia de Limarí'
$ws.Range("P20").Value = 861
$ws.Range("Q20").Value = 18

# Row 21
$ws.Range("D21").Value = 44748
$ws.Range("H21").Value = 'Cultivar IV Región'
$ws.Range("I21").Value = 'Primera'
$ws.Range("J21").Value = 250
$ws.Range("K21").Value = 17000
$ws.Range("L21").Value = 18000
$ws.Range("M21").Value = 17500
$ws.Range("N21").Value = '$/bandeja 18 kilos'
$ws.Range("O21").Value = 'Provincia de Limarí'
$ws.Range("P21").Value = 972
$ws.Range("Q21").Value = 18

# Row 22
$ws.Range("D22").Value = 44221
$ws.Range("H22").Value = 'Cultivar XV región'
$ws.Range("I22").Value = 'Primera'
$ws.Range("J22").Value = 140
$ws.Range("K22").Value = 5000
$ws.Range("L22").Value = 6000
$ws.Range("M22").Value = 5500
$ws.Range("N22").Value = '$/caja 10 kilos'
$ws.Range("O22").Value = 'Región de Arica y Parinacota'
$ws.Range("P22").Value = 550
$ws.Range("Q22").Value = 10

# Row 23
$ws.Range("D23").Value = 44454
$ws.Range("I23").Value = 'Primera'
$ws.Range("J23").Value = 160
$ws.Range("K23").Value = 19000
$ws.Range("L23").Value = 20000
$ws.Range("M23").Value = 19500
$ws.Range("P23").Value = 1083
